$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 corresponds to the "The Importance of Big Data" section, now renamed
# to "Benefits of Big Data" now that it has been completed.
$ws.Range("A8").Value = "Benefits of Big Data"

# The actual word count for that section is now known.
$ws.Range("C8").Value = 233
